$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 16 (K J O'Brien), shifting existing rows 16-23 down to 17-24
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with the new player's data
$ws.Cells.Item(16, 1).Value = "P J Moor"
$ws.Cells.Item(16, 2).Value = 4359
$ws.Cells.Item(16, 3).Value = 49
